# Update the wind-speed column (A) values in Hoja1 to reflect the new
# dataset (previous first two rows removed, remaining values shifted up).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 7
$ws.Range("A3").Value = 10
$ws.Range("A4").Value = 11.4
$ws.Range("A5").Value = 12
$ws.Range("A6").Value = 15
$ws.Range("A7").Value = 18
$ws.Range("A8").Value = 22
$ws.Range("A9").Value = 25

# The last two rows of column A no longer have values (the list got shorter)
$ws.Range("A10:A11").Clear()

# Update the active selection to match the new workbook state
$ws.Range("A2:A3").Select()
